# Fruta / hortaliza, semanal
# Insert 3 new weekly rows of Kiwi price data (Feria Lagunitas de Puerto Montt)
# at the top of the data block, pushing the existing rows (432-458) down to
# (435-461).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 432; this shifts rows 432:458 -> 435:461 and
# inherits the row-432 formatting (date number format on column D, etc.)
# for the new rows, matching native Excel "insert row" behaviour.
$ws.Rows("432:434").Insert()

# Row 432
$ws.Cells.Item(432, 1).Value = 4
$ws.Cells.Item(432, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(432, 3).Value = "Los Lagos"
$ws.Cells.Item(432, 4).Value = 45013
$ws.Cells.Item(432, 5).Value = 10
$ws.Cells.Item(432, 6).Value = "Fruta"
$ws.Cells.Item(432, 7).Value = 100101
$ws.Cells.Item(432, 8).Value = "Berries"
$ws.Cells.Item(432, 9).Value = 100101007
$ws.Cells.Item(432, 10).Value = "Kiwi"
$ws.Cells.Item(432, 11).Value = "Hayward"
$ws.Cells.Item(432, 12).Value = "Especial"
$ws.Cells.Item(432, 13).Value = 300
$ws.Cells.Item(432, 14).Value = 22000
$ws.Cells.Item(432, 15).Value = 22000
$ws.Cells.Item(432, 16).Value = 22000
$ws.Cells.Item(432, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(432, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(432, 19).Value = 1467
$ws.Cells.Item(432, 20).Value = 15

# Row 433
$ws.Cells.Item(433, 1).Value = 4
$ws.Cells.Item(433, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(433, 3).Value = "Los Lagos"
$ws.Cells.Item(433, 4).Value = 45013
$ws.Cells.Item(433, 5).Value = 10
$ws.Cells.Item(433, 6).Value = "Fruta"
$ws.Cells.Item(433, 7).Value = 100101
$ws.Cells.Item(433, 8).Value = "Berries"
$ws.Cells.Item(433, 9).Value = 100101007
$ws.Cells.Item(433, 10).Value = "Kiwi"
$ws.Cells.Item(433, 11).Value = "Hayward"
$ws.Cells.Item(433, 12).Value = "Primera"
$ws.Cells.Item(433, 13).Value = 300
$ws.Cells.Item(433, 14).Value = 19000
$ws.Cells.Item(433, 15).Value = 19000
$ws.Cells.Item(433, 16).Value = 19000
$ws.Cells.Item(433, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(433, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(433, 19).Value = 1267
$ws.Cells.Item(433, 20).Value = 15

# Row 434
$ws.Cells.Item(434, 1).Value = 4
$ws.Cells.Item(434, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(434, 3).Value = "Los Lagos"
$ws.Cells.Item(434, 4).Value = 45013
$ws.Cells.Item(434, 5).Value = 10
$ws.Cells.Item(434, 6).Value = "Fruta"
$ws.Cells.Item(434, 7).Value = 100101
$ws.Cells.Item(434, 8).Value = "Berries"
$ws.Cells.Item(434, 9).Value = 100101007
$ws.Cells.Item(434, 10).Value = "Kiwi"
$ws.Cells.Item(434, 11).Value = "Hayward"
$ws.Cells.Item(434, 12).Value = "Segunda"
$ws.Cells.Item(434, 13).Value = 300
$ws.Cells.Item(434, 14).Value = 16000
$ws.Cells.Item(434, 15).Value = 16000
$ws.Cells.Item(434, 16).Value = 16000
$ws.Cells.Item(434, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(434, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(434, 19).Value = 1067
$ws.Cells.Item(434, 20).Value = 15
